$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = "爬完cnbc從2007至今所有新聞標題及內文"

$ws.Range("A14").Value = 12.9
$ws.Range("B14").Value = "吳培瑜"
$ws.Range("C14").Value = "爬蟲 用selenium爬CNBC中market和finance相關的新聞標題"
$ws.Range("D14").Value = "code完成"
$ws.Range("E14").Value = "更快速的爬完cnbc從2007至今所有新聞標題"
$ws.Range("F14").Value = "完成書面報告與影片"

$ws.Range("E17").Select()
